$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 15).Value = "nan"
}
